$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The published "TNA BANCARIA" rate (column H) moved for a batch of plan rows
# (Pyme / PROCAMPO DIGITAL rows 48-62, and PROCAMPO rows 116-120). Columns
# J (Interes del Banco por el plazo elegido), K (Precio final financiado) and
# L (TNA con todos los gastos incluidos) are the downstream values recomputed
# from H for the same rows, so they are refreshed alongside it.
#
# Row => @(H, J, K, L)
$updates = @{
    48  = @(0.55000000000000004, 27.123287671232877,  127.75890410958904,  0.56288888888888877)
    49  = @(0.55000000000000004, 40.684931506849317,  141.38835616438354,  0.5595092592592591)
    50  = @(0.55000000000000004, 54.246575342465754,  155.01780821917805,  0.5578194444444442)
    51  = @(0.435,               21.452054794520549,  127.84426820475846,  0.56461988304093569)
    52  = @(0.46500000000000002, 34.397260273972606,  141.47080028839221,  0.56062378167641314)
    53  = @(0.48,                47.342465753424655,  155.09733237202596,  0.55862573099415214)
    54  = @(0.30499999999999999, 15.041095890410958,  127.82343987823441,  0.56419753086419777)
    55  = @(0.37,                27.36986301369863,   141.52207001522072,  0.56131687242798389)
    56  = @(0.40500000000000003, 39.945205479452056,  155.49467275494675,  0.56265432098765455)
    57  = @(0.39,                19.232876712328768,  132.22926027397261,  0.65353777777777811)
    58  = @(0.45,                33.287671232876711,  147.81602739726029,  0.64640185185185184)
    59  = @(0.495,               48.821917808219176,  165.04350684931509,  0.65946888888888922)
    60  = @(0.39,                19.232876712328768,  132.22926027397261,  0.65353777777777811)
    61  = @(0.45,                33.287671232876711,  147.81602739726029,  0.64640185185185184)
    62  = @(0.495,               48.821917808219176,  165.04350684931509,  0.65946888888888922)
    116 = @(0.56000000000000005, 9.2054794520547958,  111.38958904109589,  0.69286666666666619)
    117 = @(0.56000000000000005, 13.808219178082194,  116.08438356164383,  0.65231111111111095)
    118 = @(0.56000000000000005, 27.616438356164387,  130.16876712328767,  0.6117555555555555)
    119 = @(0.56000000000000005, 41.424657534246577,  144.25315068493151,  0.59823703703703723)
    120 = @(0.56000000000000005, 55.232876712328775,  158.33753424657536,  0.59147777777777788)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("H$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
    $ws.Range("K$row").Value = $vals[2]
    $ws.Range("L$row").Value = $vals[3]
}

# Refresh the saved selection/active-cell for the sheet view to the full
# used range, starting at A1.
$ws.Range("A1:L130").Select()
